$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the two new header cells the same formatting (bold/border/centered)
# as the existing header row before filling in their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header cells for the two new columns: I0 and IF
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2..45
$data = @{
    2  = @(7, 7)
    3  = @(9, 9)
    4  = @(7, 7)
    5  = @(7, 7)
    6  = @(8, 8)
    7  = @(3, 4)
    8  = @(5, 5)
    9  = @(9, 9)
    10 = @(6, 7)
    11 = @(8, 8)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(8, 8)
    15 = @(9, 9)
    16 = @(8, 8)
    17 = @(8, 8)
    18 = @(5, 5)
    19 = @(6, 7)
    20 = @(7, 8)
    21 = @(8, 8)
    22 = @(8, 8)
    23 = @(8, 9)
    24 = @(8, 8)
    25 = @(9, 9)
    26 = @(8, 8)
    27 = @(9, 9)
    28 = @(8, 8)
    29 = @(9, 9)
    30 = @(9, 9)
    31 = @(6, 7)
    32 = @(8, 8)
    33 = @(8, 8)
    34 = @(6, 6)
    35 = @(8, 8)
    36 = @(8, 8)
    37 = @(5, 5)
    38 = @(8, 8)
    39 = @(7, 7)
    40 = @(9, 9)
    41 = @(9, 9)
    42 = @(8, 8)
    43 = @(5, 5)
    44 = @(7, 7)
    45 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
